# Updated symbol list on Fri Dec 16 18:29:18 UTC 2022 with GitHub Actions
# Applies the refreshed coin ranking values (prices / labels / links) that
# the scraper produced for this run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2";  Value = "243.49" }
    @{ Cell = "D3";  Value = "23.48" }
    @{ Cell = "D4";  Value = "5.629" }
    @{ Cell = "D5";  Value = "0.05830" }
    @{ Cell = "D6";  Value = "3.412" }
    @{ Cell = "D7";  Value = "6.471" }
    @{ Cell = "D9";  Value = "0.7978" }
    @{ Cell = "D10"; Value = "0.1461" }
    @{ Cell = "D11"; Value = "0.07596" }
    @{ Cell = "D12"; Value = "0.03255" }
    @{ Cell = "D13"; Value = "0.02997" }
    @{ Cell = "D14"; Value = "0.09239" }
    @{ Cell = "D15"; Value = "0.001678" }
    @{ Cell = "D16"; Value = "3.324" }
    @{ Cell = "D17"; Value = "0.04744" }

    # Row 18: One -> TigerCash
    @{ Cell = "B18"; Value = "TigerCash" }
    @{ Cell = "C18"; Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch" }
    @{ Cell = "D18"; Value = "0.006246" }
    @{ Cell = "E18"; Value = "17TigerCashTCH" }

    # Row 19: TigerCash -> HotbitToken
    @{ Cell = "B19"; Value = "HotbitToken" }
    @{ Cell = "C19"; Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb" }
    @{ Cell = "D19"; Value = "0.005466" }
    @{ Cell = "E19"; Value = "18HotbitTokenHTB" }

    @{ Cell = "D20"; Value = "0.001069" }

    # Row 21: HotbitToken -> NitroEx
    @{ Cell = "B21"; Value = "NitroEx" }
    @{ Cell = "C21"; Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx" }
    @{ Cell = "D21"; Value = "0.0001501" }
    @{ Cell = "E21"; Value = "20NitroExNTX" }

    # Row 22: NitroEx -> LEO
    @{ Cell = "B22"; Value = "LEO" }
    @{ Cell = "C22"; Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo" }
    @{ Cell = "D22"; Value = "3.694" }
    @{ Cell = "E22"; Value = "21LEOLEO" }

    # Row 23: LEO -> BTSEToken
    @{ Cell = "B23"; Value = "BTSEToken" }
    @{ Cell = "C23"; Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse" }
    @{ Cell = "D23"; Value = "2.209" }
    @{ Cell = "E23"; Value = "22BTSETokenBTSE" }

    # Row 24: BTSEToken -> One
    @{ Cell = "B24"; Value = "One" }
    @{ Cell = "C24"; Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one" }
    @{ Cell = "D24"; Value = "0.01244" }
    @{ Cell = "E24"; Value = "23OneONEBestin24h" }

    @{ Cell = "D25"; Value = "0.3347" }
    @{ Cell = "D26"; Value = "0.1234" }
    @{ Cell = "D27"; Value = "0.001000" }
    @{ Cell = "E27"; Value = "26UpBotsUBXT" }

    @{ Cell = "D40"; Value = "0.04311" }
    @{ Cell = "D41"; Value = "0.007176" }

    # Row 42: BKEXToken -> CEJI
    @{ Cell = "B42"; Value = "CEJI" }
    @{ Cell = "C42"; Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji" }
    @{ Cell = "D42"; Value = "0.003603" }
    @{ Cell = "E42"; Value = "41CEJICEJI" }

    # Row 43: CEJI -> BKEXToken
    @{ Cell = "B43"; Value = "BKEXToken" }
    @{ Cell = "C43"; Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk" }
    @{ Cell = "D43"; Value = "0.1052" }
    @{ Cell = "E43"; Value = "42BKEXTokenBKK" }

    @{ Cell = "D44"; Value = "0.009536" }
    @{ Cell = "E45"; Value = "44ACDXExchangeACXTWorstin24h" }
    @{ Cell = "D46"; Value = "0.00005755" }
    @{ Cell = "D48"; Value = "0.7858" }
    @{ Cell = "D49"; Value = "0.1011" }
    @{ Cell = "E49"; Value = "48BOLOBOLO" }
    @{ Cell = "D50"; Value = "0.00002102" }
)

foreach ($update in $updates) {
    $range = $ws.Range($update.Cell)
    # Force the cell to be stored as text so that numeric-looking values
    # (e.g. "243.49", "0.001000") keep their original, exact formatting
    # instead of being coerced into a floating point number.
    $range.NumberFormat = "@"
    $range.Value = $update.Value
    # Restore the default style so the cell's appearance/format is
    # unaffected by the temporary text-number-format above.
    $range.Style = "Normal"
}
